$wb = $excel.ActiveWorkbook

# zh-cn sheet: Correspond Handoff Datetime (H2) and Correspond Handback DateTime (L2)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("H2").Value = "2016-12-16 09:18:14"
$wsZh.Range("L2").Value = "2016-12-16 09:19:06"

# de-de sheet: Correspond Handoff Datetime (H2) and Correspond Handback DateTime (L2)
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("H2").Value = "2016-12-16 09:18:29"
$wsDe.Range("L2").Value = "2016-12-16 09:19:24"

# Overview sheet: Latest HO Xliff Generate Date for de-de row (shares the same
# shared-string slot as de-de!H2 since both held the same timestamp text)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-12-16 09:18:29"
